# "update end of 2019" -> really: add a new 2020 column (J) and refresh
# several of the 2019 (column I) figures for sheet "sportKR".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New column J ("2020"). Give it the same default alignment as the
#    other year columns (C:I) before any values go in, so cells pick up
#    the existing centered style (s="1") instead of minting a new xf.
# ---------------------------------------------------------------------
$colJ = $ws.Columns.Item(10)
$colJ.HorizontalAlignment = -4108  # xlCenter

# Header + a couple of real data points that exist for 2020 already.
$ws.Range("J1").Value = 2020
$ws.Range("J2").Value = 0

# Rows 19-21 and 26-27 use the "Number" integer style (s="2") in column I;
# mirror that on the (still empty) column J cells, matching the blank
# placeholders added in those rows.
$numRows = 19, 20, 21, 26, 27
foreach ($r in $numRows) {
    $ws.Cells.Item($r, 10).NumberFormat = "0"
    $ws.Cells.Item($r, 10).HorizontalAlignment = -4108
}

# Match column J's stored width to the existing H:I columns (width 6).
$colJ.ColumnWidth = $ws.Columns.Item(9).ColumnWidth

# ---------------------------------------------------------------------
# 2) Refresh the 2019 (column I) figures that changed.
# ---------------------------------------------------------------------
$ws.Range("I9").Value = 902.2
$ws.Range("I10").Value = 48918
$ws.Range("I11").Value = 122
$ws.Range("I12").Value = 13.44
$ws.Range("I14").Value = 1
$ws.Range("I15").Value = 133.2
$ws.Range("I16").Value = 2463
$ws.Range("I17").Value = 56
$ws.Range("I18").Value = 66
$ws.Range("I25").Value = 35
$ws.Range("I28").Value = 19785

# ---------------------------------------------------------------------
# 3) Move the active selection, as in the authored workbook.
# ---------------------------------------------------------------------
$ws.Range("I26").Select()
